# Updates the cryptos list (prices / 1h volume %) to the latest scrape,
# matching the GitHub Actions commit "Updated cryptos list ... with GitHub Actions".
#
# Every touched cell holds plain text (coin names, links, price/volume strings
# such as "35.431.45" or "  +0.92%  "), never a real number. Assigning a
# numeric-looking string straight to .Value lets Excel auto-coerce it to a
# float (e.g. "244.92" -> 244.91999999999999 stored as a number), which would
# change the cell type and the stored text. Forcing NumberFormat to "@" (Text)
# before the assignment keeps it a string; ClearFormats() afterwards drops the
# now-unneeded explicit format again so the cell keeps its original (default)
# style, same as every untouched cell around it.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: "35.284.35" -> "35.431.45", E2: "  +0.48%  " -> "  +0.92%  "
Set-TextValue $ws.Range("D2") '35.431.45'
Set-TextValue $ws.Range("E2") '  +0.92%  '

# Row 3: D3: "1.898.43" -> "1.904.57", E3: "  +2.34%  " -> "  +2.64%  "
Set-TextValue $ws.Range("D3") '1.904.57'
Set-TextValue $ws.Range("E3") '  +2.64%  '

# Row 4: E4: "  -0.27%  " -> "  -0.47%  "
Set-TextValue $ws.Range("E4") '  -0.47%  '

# Row 5: D5: "244.62" -> "244.92", E5: "  +2.65%  " -> "  +2.75%  "
Set-TextValue $ws.Range("D5") '244.92'
Set-TextValue $ws.Range("E5") '  +2.75%  '

# Row 6: D6: "0.658" -> "0.661", E6: "  +5.92%  " -> "  +6.37%  "
Set-TextValue $ws.Range("D6") '0.661'
Set-TextValue $ws.Range("E6") '  +6.37%  '

# Row 7: E7: "  -0.28%  " -> "  -0.47%  "
Set-TextValue $ws.Range("E7") '  -0.47%  '

# Row 8: D8: "41.46" -> "41.78", E8: "  -1.63%  " -> "  -0.65%  "
Set-TextValue $ws.Range("D8") '41.78'
Set-TextValue $ws.Range("E8") '  -0.65%  '

# Row 9: D9: "0.351" -> "0.352", E9: "  +6.89%  " -> "  +7.28%  "
Set-TextValue $ws.Range("D9") '0.352'
Set-TextValue $ws.Range("E9") '  +7.28%  '

# Row 10: D10: "52.20" -> "52.48", E10: "  +11.65%  " -> "  +12.08%  "
Set-TextValue $ws.Range("D10") '52.48'
Set-TextValue $ws.Range("E10") '  +12.08%  '

# Row 11: D11: "0.0715" -> "0.0716", E11: "  +3.20%  " -> "  +3.32%  "
Set-TextValue $ws.Range("D11") '0.0716'
Set-TextValue $ws.Range("E11") '  +3.32%  '

# Row 12: D12: "0.0996" -> "0.0997", E12: "  +0.58%  " -> "  +0.70%  "
Set-TextValue $ws.Range("D12") '0.0997'
Set-TextValue $ws.Range("E12") '  +0.70%  '

# Row 13: D13: "2.173.65" -> "2.175.77"
Set-TextValue $ws.Range("D13") '2.175.77'

# Row 14: D14: "12.02" -> "12.09", E14: "  +4.90%  " -> "  +5.23%  "
Set-TextValue $ws.Range("D14") '12.09'
Set-TextValue $ws.Range("E14") '  +5.23%  '

# Row 15: D15: "0.697" -> "0.699", E15: "  +3.31%  " -> "  +3.39%  "
Set-TextValue $ws.Range("D15") '0.699'
Set-TextValue $ws.Range("E15") '  +3.39%  '

# Row 16: D16: "1.914.17" -> "1.901.89", E16: "  +3.13%  " -> "  +2.28%  "
Set-TextValue $ws.Range("D16") '1.901.89'
Set-TextValue $ws.Range("E16") '  +2.28%  '

# Row 17: D17: "4.86" -> "4.88", E17: "  +3.19%  " -> "  +3.65%  "
Set-TextValue $ws.Range("D17") '4.88'
Set-TextValue $ws.Range("E17") '  +3.65%  '

# Row 18: D18: "35.250.77" -> "35.370.73", E18: "  +0.46%  " -> "  +0.80%  "
Set-TextValue $ws.Range("D18") '35.370.73'
Set-TextValue $ws.Range("E18") '  +0.80%  '

# Row 19: D19: "71.46" -> "71.74", E19: "  +2.15%  " -> "  +2.61%  "
Set-TextValue $ws.Range("D19") '71.74'
Set-TextValue $ws.Range("E19") '  +2.61%  '

# Row 20: D20: "0.0₃0820" -> "0.0₃0821", E20: "  +3.43%  " -> "  +3.45%  "
Set-TextValue $ws.Range("D20") '0.0₃0821'
Set-TextValue $ws.Range("E20") '  +3.45%  '

# Row 21: D21: "239.30" -> "240.57", E21: "  -0.60%  " -> "  -0.10%  "
Set-TextValue $ws.Range("D21") '240.57'
Set-TextValue $ws.Range("E21") '  -0.10%  '

# Row 22: D22: "12.48" -> "12.54", E22: "  +2.01%  " -> "  +2.45%  "
Set-TextValue $ws.Range("D22") '12.54'
Set-TextValue $ws.Range("E22") '  +2.45%  '

# Row 23: D23: "4.80" -> "4.81", E23: "  +2.06%  " -> "  +2.27%  "
Set-TextValue $ws.Range("D23") '4.81'
Set-TextValue $ws.Range("E23") '  +2.27%  '

# Row 24: E24: "  -0.27%  " -> "  -0.29%  "
Set-TextValue $ws.Range("E24") '  -0.29%  '

# Row 25: D25: "2.40" -> "2.38", E25: "  +28.81%  " -> "  +27.94%  "
Set-TextValue $ws.Range("D25") '2.38'
Set-TextValue $ws.Range("E25") '  +27.94%  '

# Row 26: D26: "2.31" -> "2.28", E26: "  +2.20%  " -> "  +0.95%  "
Set-TextValue $ws.Range("D26") '2.28'
Set-TextValue $ws.Range("E26") '  +0.95%  '

# Row 27: D27: "170.74" -> "170.69"
Set-TextValue $ws.Range("D27") '170.69'

# Row 28: D28: "8.46" -> "8.48", E28: "  +6.48%  " -> "  +6.53%  "
Set-TextValue $ws.Range("D28") '8.48'
Set-TextValue $ws.Range("E28") '  +6.53%  '

# Row 29: D29: "18.32" -> "18.46", E29: "  +3.81%  " -> "  +4.62%  "
Set-TextValue $ws.Range("D29") '18.46'
Set-TextValue $ws.Range("E29") '  +4.62%  '

# Row 30: D30: "0.126" -> "0.127", E30: "  +1.95%  " -> "  +2.62%  "
Set-TextValue $ws.Range("D30") '0.127'
Set-TextValue $ws.Range("E30") '  +2.62%  '

# Row 31: D31: "4.15" -> "4.16", E31: "  +4.22%  " -> "  +4.26%  "
Set-TextValue $ws.Range("D31") '4.16'
Set-TextValue $ws.Range("E31") '  +4.26%  '

# Row 32: D32: "0.0565" -> "0.0566", E32: "  +1.50%  " -> "  +1.75%  "
Set-TextValue $ws.Range("D32") '0.0566'
Set-TextValue $ws.Range("E32") '  +1.75%  '

# Row 33: D33: "0.939" -> "0.938", E33: "  +11.79%  " -> "  +11.47%  "
Set-TextValue $ws.Range("D33") '0.938'
Set-TextValue $ws.Range("E33") '  +11.47%  '

# Row 34: E34: "  -0.20%  " -> "  -0.43%  "
Set-TextValue $ws.Range("E34") '  -0.43%  '

# Row 35: D35: "4.12" -> "4.13", E35: "  +3.22%  " -> "  +3.31%  "
Set-TextValue $ws.Range("D35") '4.13'
Set-TextValue $ws.Range("E35") '  +3.31%  '

# Row 36: E36: "  -3.64%  " -> "  -3.37%  "
Set-TextValue $ws.Range("E36") '  -3.37%  '

# Row 37: B37: "LidoDAOToken" -> "TrustWalletToken", C37: "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo" -> "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", D37: "2.02" -> "1.36", E37: "  -0.09%  " -> "  +4.43%  "
Set-TextValue $ws.Range("B37") 'TrustWalletToken'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D37") '1.36'
Set-TextValue $ws.Range("E37") '  +4.43%  '

# Row 38: B38: "TrustWalletToken" -> "LidoDAOToken", C38: "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" -> "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", D38: "1.35" -> "2.03", E38: "  +3.90%  " -> "  +0.05%  "
Set-TextValue $ws.Range("B38") 'LidoDAOToken'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D38") '2.03'
Set-TextValue $ws.Range("E38") '  +0.05%  '

# Row 39: E39: "  +2.17%  " -> "  +1.81%  "
Set-TextValue $ws.Range("E39") '  +1.81%  '

# Row 40: D40: "0.0210" -> "0.0211", E40: "  +4.44%  " -> "  +4.87%  "
Set-TextValue $ws.Range("D40") '0.0211'
Set-TextValue $ws.Range("E40") '  +4.87%  '

# Row 41: D41: "16.44" -> "16.41", E41: "  +9.64%  " -> "  +9.98%  "
Set-TextValue $ws.Range("D41") '16.41'
Set-TextValue $ws.Range("E41") '  +9.98%  '

# Row 42: D42: "0.0650" -> "0.0648", E42: "  +16.96%  " -> "  +16.43%  "
Set-TextValue $ws.Range("D42") '0.0648'
Set-TextValue $ws.Range("E42") '  +16.43%  '

# Row 43: D43: "89.95" -> "90.40", E43: "  -0.45%  " -> "  -0.12%  "
Set-TextValue $ws.Range("D43") '90.40'
Set-TextValue $ws.Range("E43") '  -0.12%  '

# Row 44: D44: "1.344.53" -> "1.346.28", E44: "  +0.31%  " -> "  +0.41%  "
Set-TextValue $ws.Range("D44") '1.346.28'
Set-TextValue $ws.Range("E44") '  +0.41%  '

# Row 45: D45: "2.40" -> "2.42", E45: "  +3.24%  " -> "  +4.17%  "
Set-TextValue $ws.Range("D45") '2.42'
Set-TextValue $ws.Range("E45") '  +4.17%  '

# Row 46: D46: "48.22" -> "47.86", E46: "  +38.92%  " -> "  +37.79%  "
Set-TextValue $ws.Range("D46") '47.86'
Set-TextValue $ws.Range("E46") '  +37.79%  '

# Row 47: E47: "  +1.82%  " -> "  +1.96%  "
Set-TextValue $ws.Range("E47") '  +1.96%  '

# Row 48: E48: "  -0.38%  " -> "  -0.08%  "
Set-TextValue $ws.Range("E48") '  -0.08%  '

# Row 49: D49: "6.55" -> "6.58", E49: "  -0.68%  " -> "  -0.35%  "
Set-TextValue $ws.Range("D49") '6.58'
Set-TextValue $ws.Range("E49") '  -0.35%  '

# Row 50: D50: "2.083.31" -> "2.087.76", E50: "  +2.31%  " -> "  +2.39%  "
Set-TextValue $ws.Range("D50") '2.087.76'
Set-TextValue $ws.Range("E50") '  +2.39%  '

# Row 51: B51: "Cronos" -> "Gas", C51: "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" -> "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas", D51: "0.0696" -> "11.04", E51: "  +2.28%  " -> "  -11.97%  "
Set-TextValue $ws.Range("B51") 'Gas'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
Set-TextValue $ws.Range("D51") '11.04'
Set-TextValue $ws.Range("E51") '  -11.97%  '
